# Auto-generated Excel COM-interop script applying scheduled market-data update
# to the per-class Leve profit tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 578
$ws.Range("I111").Value = 427.66666
$ws.Range("J111").Value = 1480
$ws.Range("K111").Value = 1282.99998
$ws.Range("L111").Value = 4440
$ws.Range("M111").Value = 1784.00002
$ws.Range("N111").Value = -10574

$ws.Range("H123").Value = 25450
$ws.Range("J123").Value = 25450
$ws.Range("L123").Value = 25450
$ws.Range("N123").Value = -35250

$ws.Range("H126").Value = 40780
$ws.Range("J126").Value = 40780
$ws.Range("L126").Value = 40780
$ws.Range("N126").Value = -50660

$ws.Range("H127").Value = 1313.84
$ws.Range("I127").Value = 569.0769
$ws.Range("J127").Value = 2120.6667
$ws.Range("K127").Value = 1707.2307
$ws.Range("L127").Value = 6362.000100000001
$ws.Range("M127").Value = 3252.7693
$ws.Range("N127").Value = -16282.0001

$ws.Range("H132").Value = 1066938.5
$ws.Range("I132").Value = 1116735.6
$ws.Range("J132").Value = 4600
$ws.Range("K132").Value = 3350206.8
$ws.Range("L132").Value = 13800
$ws.Range("M132").Value = -3347676.8
$ws.Range("N132").Value = -18860

$ws.Range("H137").Value = 1869.7059
$ws.Range("I137").Value = 1148.75
$ws.Range("J137").Value = 3600
$ws.Range("K137").Value = 3446.25
$ws.Range("L137").Value = 10800
$ws.Range("M137").Value = -896.25
$ws.Range("N137").Value = -15900

$ws.Range("H138").Value = 11630595
$ws.Range("I138").Value = 15385559
$ws.Range("J138").Value = 8089.048
$ws.Range("K138").Value = 46156677
$ws.Range("L138").Value = 24267.144
$ws.Range("M138").Value = -46151537
$ws.Range("N138").Value = -34547.144

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H61").Value = 2054.4443
$ws.Range("I61").Value = 1486.5454
$ws.Range("J61").Value = 2946.8572
$ws.Range("K61").Value = 1486.5454
$ws.Range("L61").Value = 2946.8572
$ws.Range("M61").Value = -1274.5454
$ws.Range("N61").Value = -3370.8572

$ws.Range("H132").Value = 1271.289
$ws.Range("I132").Value = 1028.8096
$ws.Range("K132").Value = 3086.4288
$ws.Range("M132").Value = -556.4288000000001

$ws.Range("H136").Value = 2054.4443
$ws.Range("I136").Value = 1486.5454
$ws.Range("J136").Value = 2946.8572
$ws.Range("K136").Value = 4459.6362
$ws.Range("L136").Value = 8840.571599999999
$ws.Range("M136").Value = -1909.6362
$ws.Range("N136").Value = -13940.5716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1187.1282
$ws.Range("I134").Value = 1000.0769
$ws.Range("J134").Value = 1561.2307
$ws.Range("K134").Value = 3000.2307
$ws.Range("L134").Value = 4683.6921
$ws.Range("M134").Value = -465.2307000000001
$ws.Range("N134").Value = -9753.6921

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2509.8572
$ws.Range("I31").Value = 1898.2
$ws.Range("J31").Value = 6179.8
$ws.Range("K31").Value = 1898.2
$ws.Range("L31").Value = 6179.8
$ws.Range("M31").Value = -1603.2
$ws.Range("N31").Value = -6769.8

$ws.Range("H34").Value = 2509.8572
$ws.Range("I34").Value = 1898.2
$ws.Range("J34").Value = 6179.8
$ws.Range("K34").Value = 1898.2
$ws.Range("L34").Value = 6179.8
$ws.Range("M34").Value = -1696.2
$ws.Range("N34").Value = -6583.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1140
$ws.Range("I75").Value = 600
$ws.Range("J75").Value = 1275
$ws.Range("K75").Value = 1800
$ws.Range("L75").Value = 3825
$ws.Range("M75").Value = -802
$ws.Range("N75").Value = -5821

$ws.Range("H78").Value = 1140
$ws.Range("I78").Value = 600
$ws.Range("J78").Value = 1275
$ws.Range("K78").Value = 5400
$ws.Range("L78").Value = 11475
$ws.Range("M78").Value = -408
$ws.Range("N78").Value = -21459

$ws.Range("H129").Value = 1283.3334
$ws.Range("I129").Value = 1200
$ws.Range("J129").Value = 1311.1111
$ws.Range("K129").Value = 3600
$ws.Range("L129").Value = 3933.3333
$ws.Range("M129").Value = 1400
$ws.Range("N129").Value = -13933.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2709.25
$ws.Range("I122").Value = 2585.1667
$ws.Range("K122").Value = 7755.500100000001
$ws.Range("M122").Value = -5305.500100000001

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H126").Value = 2054.3635
$ws.Range("I126").Value = 2100
$ws.Range("J126").Value = 1932.6666
$ws.Range("K126").Value = 6300
$ws.Range("L126").Value = 5797.9998
$ws.Range("M126").Value = -3830
$ws.Range("N126").Value = -10737.9998

$ws.Range("H128").Value = 40780
$ws.Range("J128").Value = 40780
$ws.Range("L128").Value = 40780
$ws.Range("N128").Value = -50740

$ws.Range("H130").Value = 40666.668
$ws.Range("J130").Value = 40666.668
$ws.Range("L130").Value = 40666.668
$ws.Range("N130").Value = -50706.668

$ws.Range("H132").Value = 9406.77
$ws.Range("I132").Value = 13036.25
$ws.Range("K132").Value = 39108.75
$ws.Range("M132").Value = -36578.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2814.625
$ws.Range("I132").Value = 1380.0526
$ws.Range("J132").Value = 4911.3076
$ws.Range("K132").Value = 4140.1578
$ws.Range("L132").Value = 14733.9228
$ws.Range("M132").Value = -1610.1578
$ws.Range("N132").Value = -19793.9228

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws.Range("H137").Value = 46749.75
$ws.Range("J137").Value = 46749.75
$ws.Range("L137").Value = 46749.75
$ws.Range("N137").Value = -56949.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H125").Value = 43886.07
$ws.Range("J125").Value = 43886.07
$ws.Range("L125").Value = 43886.07
$ws.Range("N125").Value = -53726.07

$ws.Range("H126").Value = 951.3333
$ws.Range("I126").Value = 951.3333
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 2853.9999
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -383.9998999999998
$ws.Range("N126").ClearContents()

$ws.Range("H127").Value = 45000
$ws.Range("J127").Value = 45000
$ws.Range("L127").Value = 45000
$ws.Range("N127").Value = -54920

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws.Range("H129").Value = 49214.5
$ws.Range("J129").Value = 49214.5
$ws.Range("L129").Value = 49214.5
$ws.Range("N129").Value = -59214.5

$ws.Range("H131").Value = 75000
$ws.Range("I131").Value = 70000
$ws.Range("J131").Value = 80000
$ws.Range("K131").Value = 70000
$ws.Range("L131").Value = 80000
$ws.Range("M131").Value = -64960
$ws.Range("N131").Value = -90080

$ws.Range("H136").Value = 3044
$ws.Range("I136").Value = 731.16
$ws.Range("J136").Value = 5672.227
$ws.Range("K136").Value = 2193.48
$ws.Range("L136").Value = 17016.681
$ws.Range("M136").Value = 356.52
$ws.Range("N136").Value = -22116.681

Write-Host "Applied scheduled profit-sheet update."
